# Update countries & provincias Spain
#
# Refreshes the COVID-19 "Pais" sheet to the 31-Mar-2020 09:55 data pull:
#   - headline timestamp cell (A1)
#   - updated case counts for several existing countries
#   - two newly-tracked countries ("Mauricio" and "Consejo Danes para los
#     Refugiados") inserted into the ranking, which pushes the countries
#     below them down by one row (their stats move down with them)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Headline "last updated" timestamp -------------------------------------
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 31 de Marzo de 2020 a las 09:55"

# --- Simple case-count refreshes (country/row unchanged) -------------------
# Israel (row 20): Casos activos, Recuperados
$ws.Cells.Item(20, 4).Value = 163
$ws.Cells.Item(20, 5).Value = 4650

# Noruega (row 23): Casos totales, Nuevos casos, Recuperados
$ws.Cells.Item(23, 2).Value = 4464
$ws.Cells.Item(23, 3).Value = 19
$ws.Cells.Item(23, 5).Value = 4420

# Chequia (row 25): Casos criticos
$ws.Cells.Item(25, 6).Value = 64

# Rumania (row 31): Recuperados, Muertes hoy, Muertes
$ws.Cells.Item(31, 5).Value = 1832
$ws.Cells.Item(31, 7).Value = 3
$ws.Cells.Item(31, 8).Value = 68

# Luxemburgo (row 32): Casos activos, Recuperados, Casos criticos
$ws.Cells.Item(32, 4).Value = 80
$ws.Cells.Item(32, 5).Value = 1886
$ws.Cells.Item(32, 6).Value = 31

# Finlandia (row 41): Casos totales, Nuevos casos, Recuperados
$ws.Cells.Item(41, 2).Value = 1384
$ws.Cells.Item(41, 3).Value = 32
$ws.Cells.Item(41, 5).Value = 1361

# Marruecos (row 65): Casos totales, Nuevos casos, Recuperados
$ws.Cells.Item(65, 2).Value = 574
$ws.Cells.Item(65, 3).Value = 18
$ws.Cells.Item(65, 5).Value = 526

# --- "Mauricio" inserted into the ranking around row 104 --------------------
# Honduras/Venezuela/Nigeria each shift down one row (their old stats move
# with them); Mauricio's old row further down is reused for its fresh stats.
$ws.Cells.Item(104, 1).Value = "Mauricio"
$ws.Cells.Item(104, 2).Value = 143
$ws.Cells.Item(104, 3).Value = 15
$ws.Cells.Item(104, 4).Value = 0
$ws.Cells.Item(104, 5).Value = 140
$ws.Cells.Item(104, 6).Value = 1
$ws.Cells.Item(104, 8).Value = 3

$ws.Cells.Item(105, 1).Value = "Honduras"
$ws.Cells.Item(105, 2).Value = 141
$ws.Cells.Item(105, 3).Value = 2
$ws.Cells.Item(105, 4).Value = 3
$ws.Cells.Item(105, 5).Value = 131
$ws.Cells.Item(105, 6).Value = 4
$ws.Cells.Item(105, 8).Value = 7

$ws.Cells.Item(106, 1).Value = "Venezuela"
$ws.Cells.Item(106, 2).Value = 135
$ws.Cells.Item(106, 4).Value = 39
$ws.Cells.Item(106, 5).Value = 93
$ws.Cells.Item(106, 6).Value = 6
$ws.Cells.Item(106, 8).Value = 3

$ws.Cells.Item(107, 1).Value = "Nigeria"
$ws.Cells.Item(107, 2).Value = 131
$ws.Cells.Item(107, 4).Value = 8
$ws.Cells.Item(107, 5).Value = 121
$ws.Cells.Item(107, 6).Value = 0
$ws.Cells.Item(107, 8).Value = 2

# --- "Consejo Danes para los Refugiados" inserted around row 117 -----------
# Martinica/Trinidad yTobago/Mayotte each shift down one row in the same way.
$ws.Cells.Item(117, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(117, 2).Value = 98
$ws.Cells.Item(117, 3).Value = 17
$ws.Cells.Item(117, 4).Value = 2
$ws.Cells.Item(117, 5).Value = 88
$ws.Cells.Item(117, 6).Value = 0
$ws.Cells.Item(117, 8).Value = 8

$ws.Cells.Item(118, 1).Value = "Martinica"
$ws.Cells.Item(118, 2).Value = 93
$ws.Cells.Item(118, 3).Value = 0
$ws.Cells.Item(118, 4).Value = 0
$ws.Cells.Item(118, 5).Value = 92
$ws.Cells.Item(118, 6).Value = 12
$ws.Cells.Item(118, 8).Value = 1

$ws.Cells.Item(119, 1).Value = "Trinidad yTobago"
$ws.Cells.Item(119, 2).Value = 85
$ws.Cells.Item(119, 3).Value = 2
$ws.Cells.Item(119, 4).Value = 1
$ws.Cells.Item(119, 5).Value = 81
$ws.Cells.Item(119, 6).Value = 0
$ws.Cells.Item(119, 8).Value = 3

$ws.Cells.Item(120, 1).Value = "Mayotte"
$ws.Cells.Item(120, 2).Value = 82
$ws.Cells.Item(120, 4).Value = 10
$ws.Cells.Item(120, 5).Value = 72
$ws.Cells.Item(120, 6).Value = 3
$ws.Cells.Item(120, 8).Value = 0
